# "New ideas and lets try them"
# Adds a new "60kb" test case column, renames the size headers, adds a
# "Kysimuseupload" (question/upload) test case with its own size columns,
# and records "puudu" (missing) results for the new rows/cols on Sheet9.
# Also moves the active-tab focus from Sheet4 to Sheet9.

$wb = $excel.ActiveWorkbook

# --- Sheet9: new columns / relabeled headers -----------------------------
$ws = $wb.Worksheets.Item("Sheet9")

# Existing header row (A1:E1) - relabel the size columns, add new ones
$ws.Range("C1").Value = "60kb"
$ws.Range("D1").Value = "1mb"
$ws.Range("E1").Value = "2,7mb"
$ws.Range("F1").Value = "Kysimuseupload"
$ws.Range("G1").Value = "1mb"
$ws.Range("H1").Value = "2,2mb"
$ws.Range("I1").Value = "4,8mb"

# New "Kysimuseupload" case-number column + "puudu" result columns
$ws.Range("F2").Value = 1
$ws.Range("G2:I2").Value = "puudu"

$ws.Range("F3").Value = 2
$ws.Range("F3").NumberFormat = "@"
$ws.Range("G3:I3").Value = "puudu"

$ws.Range("F4").Value = 3
$ws.Range("F4").NumberFormat = "@"
$ws.Range("G4:I4").Value = "puudu"

# Widen the new "Kysimuseupload" column
$ws.Columns.Item(6).ColumnWidth = 20.25

# --- Tab focus moves from Sheet4 to Sheet9 --------------------------------
$ws.Activate() | Out-Null
$ws.Range("I4").Select() | Out-Null

$excel.ActiveWindow.TabRatio = 0.293
